# SORMAS_User_Rights.xlsx update
# - Add 4 new user rights rows (SEE_PERSONAL_DATA_IN_JURISDICTION,
#   SEE_PERSONAL_DATA_OUTSIDE_JURISDICTION, SEE_SENSITIVE_DATA_IN_JURISDICTION,
#   SEE_SENSITIVE_DATA_OUTSIDE_JURISDICTION) to the "User Rights" sheet.
# - Bump the SORMAS Version string on the "About" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("User Rights")
$about = $wb.Worksheets.Item("About")

$xlPasteFormats = -4122

# Template cells whose formatting we reuse so the existing (shared) cell
# styles get referenced instead of new ones being minted:
#   A2 -> bold "User Right" name style
#   D2 -> green "Yes" style
#   C2 -> (also green "Yes") - kept as alt reference
#   O2 -> red "No" style
$boldTemplate = $ws.Range("A2")
$yesTemplate = $ws.Range("D2")
$noTemplate = $ws.Range("O2")

function Add-UserRightRow([int]$row, [string]$name, [string[]]$yesCols) {
    $nameCellA = $ws.Cells.Item($row, 1)
    $nameCellA.Value = $name
    $boldTemplate.Copy()
    $nameCellA.PasteSpecial($xlPasteFormats)

    $nameCellB = $ws.Cells.Item($row, 2)
    $nameCellB.Value = $name

    $allCols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y")
    foreach ($col in $allCols) {
        $cell = $ws.Range($col + $row)
        if ($yesCols -contains $col) {
            $cell.Value = "Yes"
            $yesTemplate.Copy()
            $cell.PasteSpecial($xlPasteFormats)
        } else {
            $cell.Value = "No"
            $noTemplate.Copy()
            $cell.PasteSpecial($xlPasteFormats)
        }
    }
}

$inJurisdictionYesCols = @("D","E","F","G","H","I","J","K","L","M","T","U","V")

Add-UserRightRow 108 "SEE_PERSONAL_DATA_IN_JURISDICTION" $inJurisdictionYesCols
Add-UserRightRow 109 "SEE_PERSONAL_DATA_OUTSIDE_JURISDICTION" @()
Add-UserRightRow 110 "SEE_SENSITIVE_DATA_IN_JURISDICTION" $inJurisdictionYesCols
Add-UserRightRow 111 "SEE_SENSITIVE_DATA_OUTSIDE_JURISDICTION" @()

# Bump the version shown on the About sheet.
$about.Range("A2").Value = "1.42.0-SNAPSHOT"
